$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data cells (column F: completion mark dates, column G/H: timing notes) ---
# Text cells are written in the order that reproduces the original sharedStrings.xml
# ordering (33: Просрочено..., 34: Окончательно..., 35: Продлено до ...16:00, 36: Продлено до ...22:00)

# Row 4 (# = 3) text first so "Просрочено..." / "Окончательно..." become strings 33/34
$ws.Range("G4").Value = "Просрочено на 00:55"
$ws.Range("H4").Value = "Окончательно выполнено 04.10.15 0:00"

# Row 2 (# = 1) -> string 35
$ws.Range("G2").Value = "Продлено до 04.10.15 16:00"

# Row 7 (# = 6) -> string 36
$ws.Range("G7").Value = "Продлено до 04.10.15 22:00 согласно замечаниям Влада"

# Row 2 (# = 1)
$ws.Range("F2").Value = 42280.75
$ws.Range("F2").NumberFormat = "m/d/yy h:mm"

# Row 3 (# = 2)
$ws.Range("F3").Value = 42278.895833333336
$ws.Range("F3").NumberFormat = "m/d/yy h:mm"

# Row 4 (# = 3)
$ws.Range("F4").Value = 42280.788194444445
$ws.Range("F4").NumberFormat = "m/d/yy h:mm"

# Row 7 (# = 6)
$ws.Range("F7").Value = 42280.831944444442
$ws.Range("F7").NumberFormat = "m/d/yy h:mm"

# Row 9 (# = 8)
$ws.Range("F9").Value = 42280.798611111109
$ws.Range("F9").NumberFormat = "m/d/yy h:mm"

# --- Highlight the "#" column for selected rows ---
# Green fill (FF92D050) for rows 3, 4, 8
$ws.Range("A3").Interior.Color = 5296274
$ws.Range("A4").Interior.Color = 5296274
$ws.Range("A8").Interior.Color = 5296274

# Yellow fill (FFFFFF00) for rows 2, 7, 9
$ws.Range("A2").Interior.Color = 65535
$ws.Range("A7").Interior.Color = 65535
$ws.Range("A9").Interior.Color = 65535

# --- Column widths for the new / widened columns ---
$ws.Columns.Item(7).ColumnWidth = 52
$ws.Columns.Item(8).ColumnWidth = 36.6667

# --- Update the active selection to reflect where editing left off ---
$ws.Range("H12").Select() | Out-Null
